$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1009.05
$ws.Range("I28").Value = 586.58826
$ws.Range("J28").Value = 3403
$ws.Range("K28").Value = 586.58826
$ws.Range("L28").Value = 3403
$ws.Range("M28").Value = -101.58826
$ws.Range("N28").Value = -4373

$ws.Range("H58").Value = 3103.5625
$ws.Range("I58").Value = 155.7
$ws.Range("J58").Value = 8016.6665
$ws.Range("K58").Value = 467.1
$ws.Range("L58").Value = 24049.9995
$ws.Range("M58").Value = -317.1
$ws.Range("N58").Value = -24349.9995

$ws.Range("H94").Value = 1814.75
$ws.Range("I94").Value = 1814.75
$ws.Range("K94").Value = 1814.75
$ws.Range("M94").Value = -1363.75

$ws.Range("H137").Value = 2287.5881
$ws.Range("I137").Value = 1924.4166
$ws.Range("K137").Value = 5773.2498
$ws.Range("M137").Value = -3223.2498

$ws.Range("H138").Value = 6460.6577
$ws.Range("J138").Value = 6475.25
$ws.Range("L138").Value = 19425.75
$ws.Range("N138").Value = -29705.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5735838.5
$ws.Range("I32").Value = 5735838.5
$ws.Range("K32").Value = 5735838.5
$ws.Range("M32").Value = -5735551.5

$ws.Range("H44").Value = 51377
$ws.Range("I44").Value = 20045
$ws.Range("K44").Value = 20045
$ws.Range("M44").Value = -19557

$ws.Range("H61").Value = 32266648
$ws.Range("I61").Value = 6241.75
$ws.Range("K61").Value = 6241.75
$ws.Range("M61").Value = -6029.75

$ws.Range("H102").Value = 15387813
$ws.Range("I102").Value = 22224400
$ws.Range("J102").Value = 5494.5
$ws.Range("K102").Value = 22224400
$ws.Range("L102").Value = 5494.5
$ws.Range("M102").Value = -22222778
$ws.Range("N102").Value = -8738.5

$ws.Range("H136").Value = 32266648
$ws.Range("I136").Value = 6241.75
$ws.Range("K136").Value = 18725.25
$ws.Range("M136").Value = -16175.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5214574
$ws.Range("I134").Value = 10872725
$ws.Range("K134").Value = 32618175
$ws.Range("M134").Value = -32615640

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14177.429
$ws.Range("J31").Value = 16272
$ws.Range("L31").Value = 16272
$ws.Range("N31").Value = -16862

$ws.Range("H34").Value = 14177.429
$ws.Range("J34").Value = 16272
$ws.Range("L34").Value = 16272
$ws.Range("N34").Value = -16676

$ws.Range("H41").Value = 22565.312

$ws.Range("H62").Value = 6668.3335
$ws.Range("I62").Value = 4999
$ws.Range("J62").Value = 7503
$ws.Range("K62").Value = 4999
$ws.Range("L62").Value = 7503
$ws.Range("M62").Value = -4375
$ws.Range("N62").Value = -8751

$ws.Range("H64").Value = 50000
$ws.Range("J64").Value = 50000
$ws.Range("L64").Value = 50000
$ws.Range("N64").Value = -50496

$ws.Range("H65").Value = 6668.3335
$ws.Range("I65").Value = 4999
$ws.Range("J65").Value = 7503
$ws.Range("K65").Value = 24995
$ws.Range("L65").Value = 37515
$ws.Range("M65").Value = -21875
$ws.Range("N65").Value = -43755

$ws.Range("H67").Value = 50000
$ws.Range("J67").Value = 50000
$ws.Range("L67").Value = 50000
$ws.Range("N67").Value = -51716

$ws.Range("H99").Value = 5633.6113
$ws.Range("I99").Value = 3981.6667
$ws.Range("J99").Value = 6459.5835
$ws.Range("K99").Value = 3981.6667
$ws.Range("L99").Value = 6459.5835
$ws.Range("M99").Value = -2483.6667
$ws.Range("N99").Value = -9455.583500000001

$ws.Range("H122").Value = 16669335
$ws.Range("I122").Value = 20835600
$ws.Range("K122").Value = 62506800
$ws.Range("M122").Value = -62504350

$ws.Range("H126").Value = 5633.6113
$ws.Range("I126").Value = 3981.6667
$ws.Range("J126").Value = 6459.5835
$ws.Range("K126").Value = 11945.0001
$ws.Range("L126").Value = 19378.7505
$ws.Range("M126").Value = -9475.000100000001
$ws.Range("N126").Value = -24318.7505

$ws.Range("H132").Value = 8483.147999999999
$ws.Range("I132").Value = 6617.8125
$ws.Range("K132").Value = 19853.4375
$ws.Range("M132").Value = -17323.4375

$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").Value = $null

$ws.Range("H141").Value = 100706.86
$ws.Range("I141").Value = 55000
$ws.Range("J141").Value = 108324.664
$ws.Range("K141").Value = 55000
$ws.Range("L141").Value = 108324.664
$ws.Range("M141").Value = -49820
$ws.Range("N141").Value = -118684.664

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 80572.36
$ws.Range("I2").Value = 226.86667
$ws.Range("J2").Value = 201090.6
$ws.Range("K2").Value = 1361.20002
$ws.Range("L2").Value = 1206543.6
$ws.Range("M2").Value = -1248.20002
$ws.Range("N2").Value = -1206769.6

$ws.Range("H7").Value = 261.66666
$ws.Range("I7").Value = 124.3
$ws.Range("J7").Value = 536.4
$ws.Range("K7").Value = 372.9
$ws.Range("L7").Value = 1609.2
$ws.Range("M7").Value = -260.9
$ws.Range("N7").Value = -1833.2

$ws.Range("H17").Value = 301
$ws.Range("I17").Value = 301
$ws.Range("K17").Value = 903
$ws.Range("M17").Value = -734

$ws.Range("H34").Value = 1282
$ws.Range("J34").Value = 9890
$ws.Range("L34").Value = 29670
$ws.Range("N34").Value = -29838

$ws.Range("H39").Value = 15457.333
$ws.Range("J39").Value = 17598.8
$ws.Range("L39").Value = 52796.39999999999
$ws.Range("N39").Value = -53384.39999999999

$ws.Range("H55").Value = 20009600
$ws.Range("J55").Value = 25011748
$ws.Range("L55").Value = 75035244
$ws.Range("N55").Value = -75035598

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 9666.5
$ws.Range("I46").Value = 9666.5
$ws.Range("K46").Value = 9666.5
$ws.Range("M46").Value = -9510.5

$ws.Range("H57").Value = 66666.664
$ws.Range("J57").Value = 80000
$ws.Range("L57").Value = 80000
$ws.Range("N57").Value = -81640

$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").Value = $null

$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").Value = $null

$ws.Range("H132").Value = 4885.3057
$ws.Range("I132").Value = 2953.3333
$ws.Range("K132").Value = 8859.999899999999
$ws.Range("M132").Value = -6329.999899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 799.5
$ws.Range("I55").Value = 643.8
$ws.Range("K55").Value = 643.8
$ws.Range("M55").Value = -470.8

$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").Value = $null

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 16209416
$ws.Range("I81").Value = 1168600.9
$ws.Range("J81").Value = 50051250
$ws.Range("K81").Value = 2337201.8
$ws.Range("L81").Value = 100102500
$ws.Range("M81").Value = -2336140.8
$ws.Range("N81").Value = -100104622

$ws.Range("H84").Value = 16209416
$ws.Range("I84").Value = 1168600.9
$ws.Range("J84").Value = 50051250
$ws.Range("K84").Value = 11686009
$ws.Range("L84").Value = 500512500
$ws.Range("M84").Value = -11680705
$ws.Range("N84").Value = -500523108

$ws.Range("H107").Value = 11495744
$ws.Range("I107").Value = 1128.5264
$ws.Range("K107").Value = 3385.5792
$ws.Range("M107").Value = -1465.5792

$ws.Range("H122").Value = 216579
$ws.Range("I122").Value = 312384.7
$ws.Range("K122").Value = 937154.1000000001
$ws.Range("M122").Value = -934704.1000000001

$ws.Range("H132").Value = 13355.667
$ws.Range("I132").Value = 10003.962
$ws.Range("J132").Value = 100500
$ws.Range("K132").Value = 30011.886
$ws.Range("L132").Value = 301500
$ws.Range("M132").Value = -27481.886
$ws.Range("N132").Value = -306560

$ws.Range("H136").Value = 23492772
$ws.Range("I136").Value = 43480310
$ws.Range("K136").Value = 130440930
$ws.Range("M136").Value = -130438380
